$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 0.006074256896972656
$ws.Range("D5").Value = 1066.290228157043
$ws.Range("E5").Value = 1149.390260863304
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = 259.626884
$ws.Range("H5").Value = 156.183592
$ws.Range("I5").Value = 56.49164
$ws.Range("J5").Value = "-"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 25
$ws.Range("C6").Value = 0.006862058639526367
$ws.Range("D6").Value = 1052.766982545853
$ws.Range("E6").Value = 1200.054433660507
$ws.Range("F6").Value = 1200.006325531006
$ws.Range("G6").Value = 202.577952
$ws.Range("H6").Value = 18.881952
$ws.Range("I6").Value = 26.952764
$ws.Range("J6").Value = 54.40656
